# Update cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.483.90"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "'2.553.86"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'573.42"
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("D6").Value = "'151.26"
$ws.Range("E6").Value = "  +8.89%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "'2.548.55"
$ws.Range("E9").Value = "  +4.58%  "
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "'5.75"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "'0.359"
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").Value = "'28.30"
$ws.Range("E14").Value = "  +8.68%  "
$ws.Range("D15").Value = "'3.009.24"
$ws.Range("E15").Value = "  +4.86%  "
$ws.Range("D16").Value = "'63.404.38"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "'2.561.22"
$ws.Range("E18").Value = "  +5.28%  "
$ws.Range("D19").Value = "'11.65"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("D20").Value = "'342.32"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'4.38"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").Value = "'6.89"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'66.14"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  +5.69%  "
$ws.Range("D27").Value = "'8.49"
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +9.96%  "
$ws.Range("E30").Value = "  +12.70%  "
$ws.Range("D31").Value = "'0.0₃0836"
$ws.Range("E31").Value = "  +6.08%  "
$ws.Range("E32").Value = "  +3.82%  "
$ws.Range("D33").Value = "'176.40"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("D35").Value = "'420.74"
$ws.Range("E35").Value = "  +14.32%  "
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").Value = "'19.16"
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("D38").Value = "'4.46"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'40.22"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'155.72"
$ws.Range("E43").Value = "  +6.29%  "
$ws.Range("D44").Value = "'3.82"
$ws.Range("E44").Value = "  +3.62%  "
$ws.Range("D45").Value = "'21.28"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "'0.609"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").Value = "'0.0534"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "'0.0237"
$ws.Range("E49").Value = "  +6.58%  "
$ws.Range("D50").Value = "'18.73"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").Value = "'1.83"
$ws.Range("E51").Value = "  +6.20%  "
